$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "nationalFlagHat" package to "nationalFlagHats" for all rows
#    that use it (rows 16-30, column C).
for ($r = 16; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "nationalFlagHat") {
        $cell.Value = "nationalFlagHats"
    }
}

# 2. Fix the Sri Lanka resource file names (row 28): SriLanka -> Srilanka
$ws.Range("D28").Value = "Srilanka.png"
$ws.Range("F28").Value = "Srilanka_climb.png"

# 3. Update the saved selection/active cell to D28
$ws.Range("D28").Select()
